# Applies the "Updated steps for the first and second section" edit:
#   1. "...before clicking on the "Log in" button." -> "...the "Login" button."
#   2. "After clicking on "Forgot password" button, he" ->
#        "After clicking on the "Forgot password" button, he"
#   3. "...an email with the instructions how to reset..." ->
#        "...an email with instructions how to reset..."

$d = $word.ActiveDocument

# 1) "Log in" -> "Login" (first bullet, login flow description)
$d.Content.Find.Execute(
    "the “Log in” button.", $true, $false, $false, $false, $false,
    $true, 1, $false, "the “Login” button.", 2) | Out-Null

# 2) Insert "the" before the "Forgot password" quoted button name
$d.Content.Find.Execute(
    "After clicking on “Forgot password” button, he", $true, $false, $false, $false, $false,
    $true, 1, $false, "After clicking on the “Forgot password” button, he", 2) | Out-Null

# 3) Drop the article "the" before "instructions"
$d.Content.Find.Execute(
    "an email with the instructions how to reset the password.",
    $true, $false, $false, $false, $false,
    $true, 1, $false, "an email with instructions how to reset the password.", 2) | Out-Null
